$wb = $excel.ActiveWorkbook

# Sheet 1 = "2025"
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 829.8121061896435
$ws.Range("E2").Value = 23140.84314442363
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 12936.85348562943
$ws.Range("L2").Value = 38793.7771360128
$ws.Range("M2").Value = 8710.977021824001
$ws.Range("N2").Value = 5677.290702752102
$ws.Range("O2").Value = 5600.52064569515

# Sheet 2 = "2030"
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3326.071192682715
$ws.Range("E2").Value = 36793.52723445751
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 29663.30255951012
$ws.Range("L2").Value = 43875.22642653178
$ws.Range("M2").Value = 14345.71620898827
$ws.Range("N2").Value = 7246.365569547914
$ws.Range("O2").Value = 7792.205117936748

# Sheet 3 = "2035"
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 2203.454043650656
$ws.Range("B2").Value = 5095.128507301514
$ws.Range("E2").Value = 45965.96245611054
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 41972.58945122286
$ws.Range("L2").Value = 43875.22642653178
$ws.Range("M2").Value = 18015.64306159729
$ws.Range("N2").Value = 10465.50319727667
$ws.Range("O2").Value = 10302.25342564116

# Sheet 4 = "2040"
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = 2203.454043650656
$ws.Range("B2").Value = 5095.128507301514
$ws.Range("E2").Value = 45965.96245611054
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 41972.58945122286
$ws.Range("L2").Value = 43875.22642653178
$ws.Range("M2").Value = 18015.64306159729
$ws.Range("N2").Value = 10560.45753314467
$ws.Range("O2").Value = 10302.25342564116

# Sheet 5 = "2045"
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = 4570.520850279678
$ws.Range("B2").Value = 5095.128507301514
$ws.Range("E2").Value = 45965.96245611054
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 41972.58945122286
$ws.Range("L2").Value = 43875.22642653178
$ws.Range("M2").Value = 18015.64306159729
$ws.Range("N2").Value = 10923.79731980866
$ws.Range("O2").Value = 11973.13874139941

# Sheet 6 = "2050"
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = 4570.520850279678
$ws.Range("B2").Value = 5095.128507301514
$ws.Range("E2").Value = 45965.96245611054
$ws.Range("G2").Value = 6476.740570129467
$ws.Range("I2").Value = 41972.58945122286
$ws.Range("L2").Value = 43875.22642653178
$ws.Range("M2").Value = 18015.64306159729
$ws.Range("N2").Value = 10923.79731980866
$ws.Range("O2").Value = 11973.13874139941
